$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.760.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.272.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "78.60"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.93%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.645"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.01"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("E13").Value = "  -1.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.611.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -3.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.280.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.655.74"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.26%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "233.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("E27").Value = "  -4.90%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.02%  "
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0852"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.32%  "
$ws.Range("E33").Value = "  -4.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  -6.00%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0303"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.73%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.49"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "115.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +17.93%  "
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.37"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.72%  "
$ws.Range("E45").Value = "  -3.59%  "
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.13%  "
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -4.39%  "
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("E51").Value = "  -1.79%  "
